$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Semestre ideal: EA-6 -> EA-7
$ws.Range("B9").Value2 = "EA-7"
$ws.Range("C9").Value2 = "EA-7"

# Requisitos: replace first requirement text, remove the second (whole row 26)
$novoRequisito = "LOB1217 -  Operações Unitárias e Processos  (Requisito fraco)`n"
$ws.Range("B25").Value2 = $novoRequisito
$ws.Range("C25").Value2 = $novoRequisito

# Remove the now-obsolete second requirement row entirely
$ws.Rows(26).Delete()
